# MandatoryFieldValidation: add two new worksheets (PostDetails and
# CreateSTP_Mandatory) with their reference/test data, after the existing
# CreateSTP sheet.

$wb = $excel.ActiveWorkbook
$createSTP = $wb.Worksheets.Item(1)

# --- New sheet: PostDetails -------------------------------------------------
$postDetails = $wb.Worksheets.Add($null, $createSTP)
$postDetails.Name = "PostDetails"

$postDetails.Range("A1").Value = "Title"
$postDetails.Range("B1").Value = "YourPost"

$postDetails.Range("A2").Value = "Plant1"
$postDetails.Range("A3").Value = "Plant2"
$postDetails.Range("A4").Value = "Plant3"
$postDetails.Range("A5").Value = "Plant4"
$postDetails.Range("A6").Value = "Plant5"

$postDetails.Range("B2").Value = "Details entered successfully for plant1"
$postDetails.Range("B3").Value = "Details entered successfully for plant2"
$postDetails.Range("B4").Value = "Details entered successfully for plant3"
$postDetails.Range("B5").Value = "Details entered successfully for plant4"
$postDetails.Range("B6").Value = "Details entered successfully for plant5"

$postDetails.Range("B2:B6").Select() | Out-Null

# --- New sheet: CreateSTP_Mandatory -----------------------------------------
$mandatory = $wb.Worksheets.Add($null, $postDetails)
$mandatory.Name = "CreateSTP_Mandatory"

$mandatory.Range("A1").Value = "STPName"
$mandatory.Range("B1").Value = "STPShortName"
$mandatory.Range("C1").Value = "STPDescription"
$mandatory.Range("D1").Value = "CommunityOrganizer"

$mandatory.Range("A2").Value = "Nimda1"
$mandatory.Range("B2").Value = "Nimda123"
$mandatory.Range("C2").Value = "This is a valid Mandatory Fields scenario. It is supposed give a toast message ""STP created successfully!"""
$mandatory.Range("D2").Value = "sharath sethu"

$mandatory.Range("A3").Value = "Nimda2"
$mandatory.Range("B3").Value = "Nimda223"
$mandatory.Range("C3").Value = "This is a valid Mandatory Fields scenario. It is supposed give a toast message ""STP created successfully!"""
$mandatory.Range("D3").Value = "sharath sethu"

$mandatory.Range("A4").Value = "Nimda3"
$mandatory.Range("B4").Value = "Nimda323"
$mandatory.Range("C4").Value = "TEst2"
$mandatory.Range("D4").Value = "sharath sethu"

$mandatory.Range("A5").Value = "Nimda4"
$mandatory.Range("B5").Value = "Nimda423"
$mandatory.Range("C5").Value = "TEst3"

$mandatory.Range("A6").Value = "Nimda5"
$mandatory.Range("B6").Value = "Nimda523"
$mandatory.Range("C6").Value = "Test4"

# Widen the description/organizer columns to fit their content, like the
# author did after typing the long text into column C.
$mandatory.Columns("C:D").EntireColumn.AutoFit() | Out-Null

$mandatory.Range("C7").Select() | Out-Null
